$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 6
$ws.Range("A3").Value = 6
$ws.Range("A4").Value = 6
$ws.Range("A5").Value = 21
$ws.Range("A6").Value = 21
$ws.Range("A7").Value = 21
